# Encounters changes. The summary screen now shows a maximum of 3 Encounters.

$wb = $excel.ActiveWorkbook
$wsPatients = $wb.Worksheets.Item("Patients")
$wsOrgs = $wb.Worksheets.Item("Organizations")

# --- Sheet "Patients": row 5 ("Patient 2") updates ---

# Trim the stray spaces inside the parentheses of the date-range label that is
# repeated across most of row 5 (the shared text used to read
# "( 01/05/2015 - 01/07/2016 )").
$dateRange = "(01/05/2015 - 01/07/2016)"
$wsPatients.Range("B5").Value = $dateRange
$wsPatients.Range("C5").Value = $dateRange
$wsPatients.Range("D5").Value = $dateRange
$wsPatients.Range("E5").Value = $dateRange
$wsPatients.Range("G5").Value = $dateRange
$wsPatients.Range("H5").Value = $dateRange
$wsPatients.Range("I5").Value = $dateRange
$wsPatients.Range("J5").Value = $dateRange
$wsPatients.Range("K5").Value = $dateRange
$wsPatients.Range("L5").Value = $dateRange
$wsPatients.Range("M5").Value = $dateRange
$wsPatients.Range("N5").Value = $dateRange

# The "Encounters" cell (F5) now calls out the new max-3 rule and is
# highlighted with a new orange fill.
$f5 = $wsPatients.Range("F5")
$f5.Value = "** >3 ** (01/05/2015 - 01/07/2016)"
$f5.Interior.Color = 39423

# The "Other Requirements" cell (O5) documents the new minimum encounters
# requirement.
$wsPatients.Range("O5").Value = "Min 4 encounters"

# Row 5 shrinks now that the wrapped text is shorter.
$wsPatients.Rows.Item(5).RowHeight = 22.35

# --- Active sheet / selection bookkeeping ---
# The workbook used to open on "Organizations"; it now opens on "Patients".
$wsPatients.Activate() | Out-Null
$wsPatients.Range("H5").Select() | Out-Null
